# "New API format-MA Team"
#
# The upstream fault-category taxonomy used by the API was renamed. The four
# shared labels used in column H ("Fault_cat") of the ME1_filt / ME2_filt
# sheets are replaced with more descriptive fault names:
#   Exhaust Valve       -> Exhaust Valve Fault
#   Start of Inj        -> Start of Injection Fault
#   Injection System    -> Injection System Fault
#   Combustion Blow-by  -> Blow-by in combustion chamber

$wb = $excel.ActiveWorkbook

# ME1_filt (H2:H34)
$ws1 = $wb.Worksheets.Item("ME1_filt")
$ws1.Range("H2:H7").Value   = "Exhaust Valve Fault"
$ws1.Range("H8:H25").Value  = "Start of Injection Fault"
$ws1.Range("H26:H28").Value = "Injection System Fault"
$ws1.Range("H29:H34").Value = "Blow-by in combustion chamber"

# View state tweaks recorded for this sheet: scrolled up a little and the
# selection moved from B29 to H29.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("H29").Select()

# ME2_filt (H2:H34)
$ws2 = $wb.Worksheets.Item("ME2_filt")
$ws2.Range("H2:H7")   .Value = "Blow-by in combustion chamber"
$ws2.Range("H8:H10")  .Value = "Injection System Fault"
$ws2.Range("H11:H28") .Value = "Start of Injection Fault"
$ws2.Range("H29:H34") .Value = "Exhaust Valve Fault"

# Column E was widened on this sheet (equipment codes are long strings), and
# the selection ended up on E8.
$ws2.Activate()
$ws2.Columns.Item(5).ColumnWidth = 25.7109375
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("E8").Select()

# Restore the tab that was active/visible when the file was saved.
$ws1.Activate()
